# Auto-generated edit script: updates raw market-price columns (H-N)
# across the 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# the latest scheduled-runner data pull.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = $null
$ws.Range("H33").Value = 93.86364
$ws.Range("I33").Value = 101.57895
$ws.Range("J33").Value = 45
$ws.Range("K33").Value = 101.57895
$ws.Range("L33").Value = 45
$ws.Range("M33").Value = 127.42105
$ws.Range("N33").Value = -503
$ws.Range("H40").Value = 1985.7142
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 2300
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 2300
$ws.Range("M40").Value = -1025
$ws.Range("N40").Value = -2650
$ws.Range("H98").Value = 969
$ws.Range("I98").Value = 619.1
$ws.Range("K98").Value = 619.1
$ws.Range("M98").Value = 878.9
$ws.Range("H122").Value = 969
$ws.Range("I122").Value = 619.1
$ws.Range("K122").Value = 1857.3
$ws.Range("M122").Value = 592.6999999999998
$ws.Range("H137").Value = 1199.2467
$ws.Range("I137").Value = 1068.0217
$ws.Range("J137").Value = 1393.9678
$ws.Range("K137").Value = 3204.0651
$ws.Range("L137").Value = 4181.903399999999
$ws.Range("M137").Value = -654.0650999999998
$ws.Range("N137").Value = -9281.903399999999
$ws.Range("H138").Value = 1353.54
$ws.Range("J138").Value = 1880.2203
$ws.Range("L138").Value = 5640.6609
$ws.Range("N138").Value = -15920.6609
$ws.Range("H141").Value = 3079.814
$ws.Range("I141").Value = 1036.4054
$ws.Range("J141").Value = 15680.833
$ws.Range("K141").Value = 3109.2162
$ws.Range("L141").Value = 47042.499
$ws.Range("M141").Value = 2070.7838
$ws.Range("N141").Value = -57402.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1320.1305
$ws.Range("I2").Value = 1378.421
$ws.Range("J2").Value = 1043.25
$ws.Range("K2").Value = 1378.421
$ws.Range("L2").Value = 1043.25
$ws.Range("M2").Value = -1265.421
$ws.Range("N2").Value = -1269.25
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = $null
$ws.Range("H32").Value = 709695.5
$ws.Range("I32").Value = 828319.8
$ws.Range("J32").Value = 16199.77
$ws.Range("K32").Value = 828319.8
$ws.Range("L32").Value = 16199.77
$ws.Range("M32").Value = -828032.8
$ws.Range("N32").Value = -16773.77
$ws.Range("H45").Value = 2468.5264
$ws.Range("I45").Value = 2766.4443
$ws.Range("K45").Value = 2766.4443
$ws.Range("M45").Value = -2389.4443
$ws.Range("H74").Value = 1986.7778
$ws.Range("I74").Value = 1110.2
$ws.Range("K74").Value = 1110.2
$ws.Range("M74").Value = -236.2
$ws.Range("H77").Value = 1986.7778
$ws.Range("I77").Value = 1110.2
$ws.Range("K77").Value = 5551
$ws.Range("M77").Value = -1183
$ws.Range("H110").Value = 45938.9
$ws.Range("I110").Value = 57217.375
$ws.Range("J110").Value = 825
$ws.Range("K110").Value = 57217.375
$ws.Range("L110").Value = 825
$ws.Range("M110").Value = -55172.375
$ws.Range("N110").Value = -4915
$ws.Range("H116").Value = 1320.1305
$ws.Range("I116").Value = 1378.421
$ws.Range("J116").Value = 1043.25
$ws.Range("K116").Value = 1378.421
$ws.Range("L116").Value = 1043.25
$ws.Range("M116").Value = 915.579
$ws.Range("N116").Value = -5631.25
$ws.Range("H121").Value = 59980
$ws.Range("J121").Value = 59980
$ws.Range("L121").Value = 59980
$ws.Range("N121").Value = -63474

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1320.1305
$ws.Range("I3").Value = 1378.421
$ws.Range("J3").Value = 1043.25
$ws.Range("K3").Value = 1378.421
$ws.Range("L3").Value = 1043.25
$ws.Range("M3").Value = -1264.421
$ws.Range("N3").Value = -1271.25
$ws.Range("H15").Value = 25000004
$ws.Range("J15").Value = 7
$ws.Range("L15").Value = 7
$ws.Range("N15").Value = -461
$ws.Range("H86").Value = 2123.7144
$ws.Range("I86").Value = 2438.353
$ws.Range("J86").Value = 1637.4546
$ws.Range("K86").Value = 2438.353
$ws.Range("L86").Value = 1637.4546
$ws.Range("M86").Value = -1315.353
$ws.Range("N86").Value = -3883.4546
$ws.Range("H89").Value = 2123.7144
$ws.Range("I89").Value = 2438.353
$ws.Range("J89").Value = 1637.4546
$ws.Range("K89").Value = 12191.765
$ws.Range("L89").Value = 8187.273
$ws.Range("M89").Value = -6575.764999999999
$ws.Range("N89").Value = -19419.273
$ws.Range("H99").Value = 968.6875
$ws.Range("I99").Value = 824.9167
$ws.Range("K99").Value = 824.9167
$ws.Range("M99").Value = 673.0833
$ws.Range("H134").Value = 2462.1333
$ws.Range("I134").Value = 2352.6428
$ws.Range("J134").Value = 3995
$ws.Range("K134").Value = 7057.928400000001
$ws.Range("L134").Value = 11985
$ws.Range("M134").Value = -4522.928400000001
$ws.Range("N134").Value = -17055

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3567.5308
$ws.Range("I31").Value = 1108.6171
$ws.Range("J31").Value = 6966.6177
$ws.Range("K31").Value = 1108.6171
$ws.Range("L31").Value = 6966.6177
$ws.Range("M31").Value = -813.6170999999999
$ws.Range("N31").Value = -7556.6177
$ws.Range("H34").Value = 3567.5308
$ws.Range("I34").Value = 1108.6171
$ws.Range("J34").Value = 6966.6177
$ws.Range("K34").Value = 1108.6171
$ws.Range("L34").Value = 6966.6177
$ws.Range("M34").Value = -906.6170999999999
$ws.Range("N34").Value = -7370.6177
$ws.Range("H58").Value = 1314.4894
$ws.Range("I58").Value = 1010.2414
$ws.Range("J58").Value = 1804.6666
$ws.Range("K58").Value = 1010.2414
$ws.Range("L58").Value = 1804.6666
$ws.Range("M58").Value = -807.2414
$ws.Range("N58").Value = -2210.6666
$ws.Range("H132").Value = 3269451.8
$ws.Range("I132").Value = 1282.6757
$ws.Range("K132").Value = 3848.0271
$ws.Range("M132").Value = -1318.0271
$ws.Range("H136").Value = 1314.4894
$ws.Range("I136").Value = 1010.2414
$ws.Range("J136").Value = 1804.6666
$ws.Range("K136").Value = 3030.7242
$ws.Range("L136").Value = 5413.9998
$ws.Range("M136").Value = -480.7242000000001
$ws.Range("N136").Value = -10513.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1423.0333
$ws.Range("I5").Value = 464.17648
$ws.Range("J5").Value = 2676.923
$ws.Range("K5").Value = 1392.52944
$ws.Range("L5").Value = 8030.768999999999
$ws.Range("M5").Value = -1280.52944
$ws.Range("N5").Value = -8254.769
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 350
$ws.Range("J17").Value = 800
$ws.Range("K17").Value = 1050
$ws.Range("L17").Value = 2400
$ws.Range("M17").Value = -881
$ws.Range("N17").Value = -2738
$ws.Range("H122").Value = 2443.111
$ws.Range("I122").Value = 377.19354
$ws.Range("J122").Value = 5227.609
$ws.Range("K122").Value = 3394.74186
$ws.Range("L122").Value = 47048.481
$ws.Range("M122").Value = -944.7418600000001
$ws.Range("N122").Value = -51948.481
$ws.Range("H135").Value = 1423.0333
$ws.Range("I135").Value = 464.17648
$ws.Range("J135").Value = 2676.923
$ws.Range("K135").Value = 4177.58832
$ws.Range("L135").Value = 24092.307
$ws.Range("M135").Value = -1642.58832
$ws.Range("N135").Value = -29162.307
$ws.Range("H137").Value = 6951526
$ws.Range("I137").Value = 23825174
$ws.Range("J137").Value = 3553.7646
$ws.Range("K137").Value = 71475522
$ws.Range("L137").Value = 10661.2938
$ws.Range("M137").Value = -71470422
$ws.Range("N137").Value = -20861.2938

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 52584.5
$ws.Range("I19").Value = 320
$ws.Range("J19").Value = 70006
$ws.Range("K19").Value = 320
$ws.Range("L19").Value = 70006
$ws.Range("M19").Value = -32
$ws.Range("N19").Value = -70582
$ws.Range("H109").Value = 10285
$ws.Range("J109").Value = 10285
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365
$ws.Range("H113").Value = 5334.76
$ws.Range("I113").Value = 1470.625
$ws.Range("J113").Value = 12204.333
$ws.Range("K113").Value = 1470.625
$ws.Range("L113").Value = 12204.333
$ws.Range("M113").Value = 699.375
$ws.Range("N113").Value = -16544.333
$ws.Range("H122").Value = 1440.3334
$ws.Range("J122").Value = 2502.6667
$ws.Range("L122").Value = 7508.000100000001
$ws.Range("N122").Value = -12408.0001
$ws.Range("H126").Value = 1168.4286
$ws.Range("I126").Value = 717
$ws.Range("J126").Value = 1507
$ws.Range("K126").Value = 2151
$ws.Range("L126").Value = 4521
$ws.Range("M126").Value = 319
$ws.Range("N126").Value = -9461
$ws.Range("H132").Value = 1779.2188
$ws.Range("I132").Value = 1419.8701
$ws.Range("J132").Value = 3235.5264
$ws.Range("K132").Value = 4259.6103
$ws.Range("L132").Value = 9706.5792
$ws.Range("M132").Value = -1729.6103
$ws.Range("N132").Value = -14766.5792

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = $null
$ws.Range("N19").Value = $null
$ws.Range("H122").Value = 4556.207
$ws.Range("H132").Value = 2715.61
$ws.Range("I132").Value = 2334.5227
$ws.Range("J132").Value = 3833.4666
$ws.Range("K132").Value = 7003.5681
$ws.Range("L132").Value = 11500.3998
$ws.Range("M132").Value = -4473.5681
$ws.Range("N132").Value = -16560.3998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4419.8887
$ws.Range("I96").Value = 3176
$ws.Range("J96").Value = 5974.75
$ws.Range("K96").Value = 3176
$ws.Range("L96").Value = 5974.75
$ws.Range("M96").Value = -1803
$ws.Range("N96").Value = -8720.75
$ws.Range("H122").Value = 2843.3333
$ws.Range("I122").Value = 2371.5
$ws.Range("J122").Value = 4191.4287
$ws.Range("K122").Value = 7114.5
$ws.Range("L122").Value = 12574.2861
$ws.Range("M122").Value = -4664.5
$ws.Range("N122").Value = -17474.2861
$ws.Range("H132").Value = 4631143.5
$ws.Range("I132").Value = 1534.1578
$ws.Range("K132").Value = 4602.4734
$ws.Range("M132").Value = -2072.4734
